# datasetOverview.xlsx — add 8 new summary columns (AO:AV) to Sheet1:
#   GenPeak, GenPeakEl, mostDivGenPeak, mostDivGenPeakEl,
#   SFPeak, SFPeakEl, mostDivSFPeak, mostDivSFPeakEl
# and update the workbook/sheet view state to reflect the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("AO", "AP", "AQ", "AR", "AS", "AT", "AU", "AV")

# --- header row (row 1) ---------------------------------------------------
$headers = @(
    "GenPeak",
    "GenPeakEl",
    "mostDivGenPeak",
    "mostDivGenPeakEl",
    "SFPeak",
    "SFPeakEl",
    "mostDivSFPeak",
    "mostDivSFPeakEl"
)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "1").Value = $headers[$i]
}

# --- data rows --------------------------------------------------------------
# Not every record row got peak-summary data (rows 7, 12, 19, 20 are left
# blank, same as the other #N/A rows in this dataset).
$rowData = @{
    2  = @(28, 1800, 39, 1800, 6, 1800, 65, 1800)
    3  = @(22, 1300, 30, 1300, 4, 800, 53, 1300)
    4  = @(27, 700, 19, 600, 5, 100, 50, 600)
    5  = @(23, 1400, 8, 1600, 6, 1400, 32, 1600)
    6  = @(23, 1400, 31, 2100, 4, 1000, 54, 1900)
    8  = @(13, 250, 8, 475, 4, 250, 20, 475)
    9  = @(26, 425, 21, 875, 8, 425, 71, 825)
    10 = @(12, 800, 7, 800, 6, 800, 13, 800)
    11 = @(16, 800, 14, 1000, 4, 600, 25, 1000)
    13 = @(14, 400, 14, 900, 4, 400, 22, 600)
    14 = @(16, 440, 6, 719, 6, 511, 16, 719)
    15 = @(17, 1000, 3, 800, 5, 1000, 12, 1000)
    16 = @(49, 300, 10, 300, 10, 300, 47, 300)
    17 = @(24, 100, 9, 100, 5, 200, 36, 100)
    18 = @(61, 300, 56, 500, 12, 50, 171, 500)
    21 = @(39, 300, 16, 300, 7, 200, 71, 300)
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $vals[$i]
    }
}

# --- window / selection state ----------------------------------------------
# Scroll the view over towards the newly-added columns and select the block
# of new peak-summary values that were just filled in.
$ws.Range("AO8:AV11").Select()

$win = $excel.ActiveWindow
$win.ScrollColumn = 35   # AI
$win.ScrollRow = 1

$excel.Left = 1640
$excel.Top = 0
